# Update rows 2-9 with recalculated TPM-based NATMI values, and add two new rows
# (FAPs -> Resolving-Mac, MuSCs -> Resolving-Mac) reflecting the new 'Resolving-Mac'
# target cluster added to the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 'FAPs'
$ws.Cells.Item(2,2).Value = 'Lgi1'
$ws.Cells.Item(2,3).Value = 'Adam23'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.014884
$ws.Cells.Item(2,8).Value = 0.044652
$ws.Cells.Item(2,9).Value = 0.6275843652054141
$ws.Cells.Item(2,10).Value = 0.6275843652054141
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.2320676666666667
$ws.Cells.Item(2,14).Value = 0.696203
$ws.Cells.Item(2,15).Value = 0.01170834749781651
$ws.Cells.Item(2,16).Value = 0.01170834749781651
$ws.Cells.Item(2,17).Value = 0.003454095150666667
$ws.Cells.Item(2,18).Value = 0.031086856356
$ws.Cells.Item(2,19).Value = 0.007347975832021573
$ws.Cells.Item(2,20).Value = 0.007347975832021573

# Row 3
$ws.Cells.Item(3,1).Value = 'FAPs'
$ws.Cells.Item(3,2).Value = 'Lgi1'
$ws.Cells.Item(3,3).Value = 'Adam23'
$ws.Cells.Item(3,4).Value = 'FAPs'
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.014884
$ws.Cells.Item(3,8).Value = 0.044652
$ws.Cells.Item(3,9).Value = 0.6275843652054141
$ws.Cells.Item(3,10).Value = 0.6275843652054141
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 17.63262733333333
$ws.Cells.Item(3,14).Value = 52.897882
$ws.Cells.Item(3,15).Value = 0.889606600882922
$ws.Cells.Item(3,16).Value = 0.8896066008829221
$ws.Cells.Item(3,17).Value = 0.2624440252293333
$ws.Cells.Item(3,18).Value = 2.361996227064
$ws.Cells.Item(3,19).Value = 0.5583031938976547
$ws.Cells.Item(3,20).Value = 0.5583031938976548

# Row 4
$ws.Cells.Item(4,1).Value = 'FAPs'
$ws.Cells.Item(4,2).Value = 'Lgi1'
$ws.Cells.Item(4,3).Value = 'Adam23'
$ws.Cells.Item(4,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.014884
$ws.Cells.Item(4,8).Value = 0.044652
$ws.Cells.Item(4,9).Value = 0.6275843652054141
$ws.Cells.Item(4,10).Value = 0.6275843652054141
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.7378426666666668
$ws.Cells.Item(4,14).Value = 2.213528
$ws.Cells.Item(4,15).Value = 0.03722585944063267
$ws.Cells.Item(4,16).Value = 0.03722585944063268
$ws.Cells.Item(4,17).Value = 0.01098205025066667
$ws.Cells.Item(4,18).Value = 0.09883845225600001
$ws.Cells.Item(4,19).Value = 0.02336236736627543
$ws.Cells.Item(4,20).Value = 0.02336236736627543

# Row 5
$ws.Cells.Item(5,1).Value = 'FAPs'
$ws.Cells.Item(5,2).Value = 'Lgi1'
$ws.Cells.Item(5,3).Value = 'Adam23'
$ws.Cells.Item(5,4).Value = 'MuSCs'
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.014884
$ws.Cells.Item(5,8).Value = 0.044652
$ws.Cells.Item(5,9).Value = 0.6275843652054141
$ws.Cells.Item(5,10).Value = 0.6275843652054141
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.143822333333333
$ws.Cells.Item(5,14).Value = 3.431467
$ws.Cells.Item(5,15).Value = 0.05770846730521116
$ws.Cells.Item(5,16).Value = 0.05770846730521117
$ws.Cells.Item(5,17).Value = 0.01702465160933333
$ws.Cells.Item(5,18).Value = 0.153221864484
$ws.Cells.Item(5,19).Value = 0.03621693182071834
$ws.Cells.Item(5,20).Value = 0.03621693182071834

# Row 6
$ws.Cells.Item(6,1).Value = 'FAPs'
$ws.Cells.Item(6,2).Value = 'Lgi1'
$ws.Cells.Item(6,3).Value = 'Adam23'
$ws.Cells.Item(6,4).Value = 'Resolving-Mac'
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.014884
$ws.Cells.Item(6,8).Value = 0.044652
$ws.Cells.Item(6,9).Value = 0.6275843652054141
$ws.Cells.Item(6,10).Value = 0.6275843652054141
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.07434200000000001
$ws.Cells.Item(6,14).Value = 0.223026
$ws.Cells.Item(6,15).Value = 0.003750724873417703
$ws.Cells.Item(6,16).Value = 0.003750724873417703
$ws.Cells.Item(6,17).Value = 0.001106506328
$ws.Cells.Item(6,18).Value = 0.009958556952
$ws.Cells.Item(6,19).Value = 0.002353896288744006
$ws.Cells.Item(6,20).Value = 0.002353896288744006

# Row 7
$ws.Cells.Item(7,1).Value = 'MuSCs'
$ws.Cells.Item(7,2).Value = 'Lgi1'
$ws.Cells.Item(7,3).Value = 'Adam23'
$ws.Cells.Item(7,4).Value = 'ECs'
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.008832333333333333
$ws.Cells.Item(7,8).Value = 0.026497
$ws.Cells.Item(7,9).Value = 0.372415634794586
$ws.Cells.Item(7,10).Value = 0.3724156347945861
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.2320676666666667
$ws.Cells.Item(7,14).Value = 0.696203
$ws.Cells.Item(7,15).Value = 0.01170834749781651
$ws.Cells.Item(7,16).Value = 0.01170834749781651
$ws.Cells.Item(7,17).Value = 0.002049698987888889
$ws.Cells.Item(7,18).Value = 0.018447290891
$ws.Cells.Item(7,19).Value = 0.004360371665794938
$ws.Cells.Item(7,20).Value = 0.004360371665794939

# Row 8
$ws.Cells.Item(8,1).Value = 'MuSCs'
$ws.Cells.Item(8,2).Value = 'Lgi1'
$ws.Cells.Item(8,3).Value = 'Adam23'
$ws.Cells.Item(8,4).Value = 'FAPs'
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.008832333333333333
$ws.Cells.Item(8,8).Value = 0.026497
$ws.Cells.Item(8,9).Value = 0.372415634794586
$ws.Cells.Item(8,10).Value = 0.3724156347945861
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 17.63262733333333
$ws.Cells.Item(8,14).Value = 52.897882
$ws.Cells.Item(8,15).Value = 0.889606600882922
$ws.Cells.Item(8,16).Value = 0.8896066008829221
$ws.Cells.Item(8,17).Value = 0.1557372421504444
$ws.Cells.Item(8,18).Value = 1.401635179354
$ws.Cells.Item(8,19).Value = 0.3313034069852673
$ws.Cells.Item(8,20).Value = 0.3313034069852674

# Row 9
$ws.Cells.Item(9,1).Value = 'MuSCs'
$ws.Cells.Item(9,2).Value = 'Lgi1'
$ws.Cells.Item(9,3).Value = 'Adam23'
$ws.Cells.Item(9,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.008832333333333333
$ws.Cells.Item(9,8).Value = 0.026497
$ws.Cells.Item(9,9).Value = 0.372415634794586
$ws.Cells.Item(9,10).Value = 0.3724156347945861
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.7378426666666668
$ws.Cells.Item(9,14).Value = 2.213528
$ws.Cells.Item(9,15).Value = 0.03722585944063267
$ws.Cells.Item(9,16).Value = 0.03722585944063268
$ws.Cells.Item(9,17).Value = 0.006516872379555556
$ws.Cells.Item(9,18).Value = 0.058651851416
$ws.Cells.Item(9,19).Value = 0.01386349207435725
$ws.Cells.Item(9,20).Value = 0.01386349207435725

# Row 10
$ws.Cells.Item(10,1).Value = 'MuSCs'
$ws.Cells.Item(10,2).Value = 'Lgi1'
$ws.Cells.Item(10,3).Value = 'Adam23'
$ws.Cells.Item(10,4).Value = 'MuSCs'
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.008832333333333333
$ws.Cells.Item(10,8).Value = 0.026497
$ws.Cells.Item(10,9).Value = 0.372415634794586
$ws.Cells.Item(10,10).Value = 0.3724156347945861
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.143822333333333
$ws.Cells.Item(10,14).Value = 3.431467
$ws.Cells.Item(10,15).Value = 0.05770846730521116
$ws.Cells.Item(10,16).Value = 0.05770846730521117
$ws.Cells.Item(10,17).Value = 0.01010262012211111
$ws.Cells.Item(10,18).Value = 0.09092358109900001
$ws.Cells.Item(10,19).Value = 0.02149153548449283
$ws.Cells.Item(10,20).Value = 0.02149153548449283

# Row 11
$ws.Cells.Item(11,1).Value = 'MuSCs'
$ws.Cells.Item(11,2).Value = 'Lgi1'
$ws.Cells.Item(11,3).Value = 'Adam23'
$ws.Cells.Item(11,4).Value = 'Resolving-Mac'
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.008832333333333333
$ws.Cells.Item(11,8).Value = 0.026497
$ws.Cells.Item(11,9).Value = 0.372415634794586
$ws.Cells.Item(11,10).Value = 0.3724156347945861
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.07434200000000001
$ws.Cells.Item(11,14).Value = 0.223026
$ws.Cells.Item(11,15).Value = 0.003750724873417703
$ws.Cells.Item(11,16).Value = 0.003750724873417703
$ws.Cells.Item(11,17).Value = 0.0006566133246666667
$ws.Cells.Item(11,18).Value = 0.005909519922
$ws.Cells.Item(11,19).Value = 0.001396828584673697
$ws.Cells.Item(11,20).Value = 0.001396828584673697
